# feat: add localization for chest and potion bonus pickup
#
# Inserts a new localization row (row 39) for the "tmp noMoreUpgrades"
# key ("No more upgrades…" and its FR/ZH/JP/ES translations), shifting
# the existing rows 39-51 down to 40-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 39 (shifts old rows 39..51 down to 40..52).
$ws.Rows.Item(39).Insert()

# Copy the formatting (style s="4", border, wrap text, vertical align)
# from the row that is now directly below (old row 39, now row 40) so the
# new row matches the rest of the table instead of getting the bare
# default column style.
$ws.Range("A40:G40").Copy()
$ws.Range("A39:G39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values.
# Column A reuses the existing "UI player stats" category string.
$ws.Range("A39").Value = "UI player stats"
$ws.Range("B39").Value = "tmp noMoreUpgrades"
$ws.Range("C39").Value = "No more upgrades…"
$ws.Range("D39").Value = "Aucun bonus restant…"
# NOTE: write G, then F, then E so the new shared strings land in the
# same order as the source workbook (uniqueCount indices 324/325/326).
$ws.Range("G39").Value = "没有更多的升级"
$ws.Range("F39").Value = "もうアップグレードはありません"
$ws.Range("E39").Value = "No más actualizaciones"

# The new row needs a taller height to fit the wrapped translations.
$ws.Rows.Item(39).RowHeight = 43.2

# Update the view to reflect where the author was looking/selecting
# after the edit.
$ws.Range("H39").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
